$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.00" or "7.77"
# are not auto-coerced into numbers by Excel's smart-entry, matching the
# original inlineStr ("text") cell type used throughout the sheet.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '60.750.24'
$ws.Range('E2').Value = '  -1.90%  '
$ws.Range('D3').Value = '2.902.38'
$ws.Range('E3').Value = '  -3.10%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '524.78'
$ws.Range('E5').Value = '  -3.51%  '
$ws.Range('D6').Value = '143.45'
$ws.Range('E6').Value = '  -6.02%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '0.543'
$ws.Range('E8').Value = '  -4.93%  '
$ws.Range('D9').Value = '2.908.75'
$ws.Range('E9').Value = '  -3.32%  '
$ws.Range('E10').Value = '  -6.10%  '
$ws.Range('D11').Value = '6.15'
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('D12').Value = '0.357'
$ws.Range('E12').Value = '  -3.47%  '
$ws.Range('D13').Value = '3.408.44'
$ws.Range('E13').Value = '  -3.19%  '
$ws.Range('E14').Value = '  +2.34%  '
$ws.Range('D15').Value = '60.779.97'
$ws.Range('E15').Value = '  -1.89%  '
$ws.Range('D16').Value = '22.44'
$ws.Range('E16').Value = '  -6.77%  '
$ws.Range('D17').Value = '2.916.53'
$ws.Range('E17').Value = '  -2.83%  '
$ws.Range('E18').Value = '  -4.79%  '
$ws.Range('D19').Value = '4.88'
$ws.Range('E19').Value = '  -5.82%  '
$ws.Range('D20').Value = '11.49'
$ws.Range('E20').Value = '  -4.91%  '
$ws.Range('D21').Value = '352.06'
$ws.Range('E21').Value = '  -7.22%  '
$ws.Range('D22').Value = '6.49'
$ws.Range('E22').Value = '  -3.25%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.26%  '
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('D25').Value = '64.76'
$ws.Range('E25').Value = '  -2.13%  '
$ws.Range('D26').Value = '0.448'
$ws.Range('E26').Value = '  -4.69%  '
$ws.Range('E27').Value = '  -5.52%  '
$ws.Range('D28').Value = '0.996'
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '7.77'
$ws.Range('E29').Value = '  -6.12%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0861'
$ws.Range('E30').Value = '  -8.31%  '
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('E32').Value = '  -3.08%  '
$ws.Range('D33').Value = '19.57'
$ws.Range('E33').Value = '  -4.62%  '
$ws.Range('D34').Value = '153.59'
$ws.Range('E34').Value = '  -4.43%  '
$ws.Range('D35').Value = '4.38'
$ws.Range('E35').Value = '  -5.35%  '
$ws.Range('D36').Value = '5.57'
$ws.Range('E36').Value = '  -7.12%  '
$ws.Range('D37').Value = '0.991'
$ws.Range('E37').Value = '  -7.73%  '
$ws.Range('D38').Value = '1.19'
$ws.Range('E38').Value = '  -7.02%  '
$ws.Range('D39').Value = '37.54'
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('D40').Value = '1.46'
$ws.Range('E40').Value = '  -6.27%  '
$ws.Range('D41').Value = '0.651'
$ws.Range('E41').Value = '  -3.60%  '
$ws.Range('D42').Value = '2.285.96'
$ws.Range('E42').Value = '  -5.52%  '
$ws.Range('D43').Value = '3.67'
$ws.Range('E43').Value = '  -6.15%  '
$ws.Range('D44').Value = '0.0581'
$ws.Range('E44').Value = '  -2.06%  '
$ws.Range('D45').Value = '20.34'
$ws.Range('E45').Value = '  -8.11%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').Value = '4.91'
$ws.Range('E47').Value = '  -5.70%  '
$ws.Range('D48').Value = '0.0236'
$ws.Range('E48').Value = '  -3.78%  '
$ws.Range('E49').Value = '  -0.70%  '
$ws.Range('D50').Value = '0.0909'
$ws.Range('E50').Value = '  -4.74%  '
$ws.Range('D51').Value = '18.37'
$ws.Range('E51').Value = '  -7.01%  '

# Restore default (Normal) style on column D so no stray number format
# is left applied to the cells (matches original unstyled cells).
$ws.Range("D2:D51").Style = "Normal"
